{"js": "// Replace the date line and each \"AAA\u00d7B=\" problem text with its new value.\n// The mapping is a 1:1, order-independent set of literal text replacements\n// (every old value in the document is unique), so we simply search for\n// each exact old string and replace it with the corresponding new string.\nconst replacements = [\n  [\"2024-04-25 Thursday\", \"2024-04-26 Friday\"],\n  [\"217\u00d78=\", \"729\u00d72=\"],\n  [\"741\u00d78=\", \"269\u00d77=\"],\n  [\"682\u00d78=\", \"569\u00d79=\"],\n  [\"485\u00d79=\", \"504\u00d76=\"],\n  [\"859\u00d73=\", \"646\u00d74=\"],\n  [\"343\u00d78=\", \"238\u00d76=\"],\n  [\"621\u00d73=\", \"581\u00d78=\"],\n  [\"569\u00d74=\", \"871\u00d72=\"],\n  [\"781\u00d78=\", \"821\u00d72=\"],\n  [\"582\u00d78=\", \"944\u00d75=\"],\n  [\"193\u00d78=\", \"643\u00d78=\"],\n  [\"277\u00d75=\", \"815\u00d77=\"],\n  [\"872\u00d76=\", \"987\u00d79=\"],\n  [\"419\u00d72=\", \"438\u00d75=\"],\n  [\"495\u00d76=\", \"262\u00d74=\"],\n  [\"158\u00d72=\", \"441\u00d74=\"],\n  [\"152\u00d77=\", \"394\u00d74=\"],\n  [\"473\u00d79=\", \"997\u00d79=\"],\n  [\"887\u00d72=\", \"790\u00d76=\"],\n  [\"400\u00d75=\", \"757\u00d77=\"],\n  [\"586\u00d73=\", \"681\u00d78=\"],\n  [\"645\u00d72=\", \"690\u00d73=\"],\n  [\"584\u00d73=\", \"321\u00d75=\"],\n  [\"813\u00d79=\", \"822\u00d76=\"],\n  [\"106\u00d77=\", \"122\u00d72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1:1 literal text replacements (every \"old\" value below is unique in the\n# document, so a plain Find/Replace on each pair is unambiguous).\n$pairs = @(\n    @(\"2024-04-25 Thursday\", \"2024-04-26 Friday\"),\n    @(\"217\u00d78=\", \"729\u00d72=\"),\n    @(\"741\u00d78=\", \"269\u00d77=\"),\n    @(\"682\u00d78=\", \"569\u00d79=\"),\n    @(\"485\u00d79=\", \"504\u00d76=\"),\n    @(\"859\u00d73=\", \"646\u00d74=\"),\n    @(\"343\u00d78=\", \"238\u00d76=\"),\n    @(\"621\u00d73=\", \"581\u00d78=\"),\n    @(\"569\u00d74=\", \"871\u00d72=\"),\n    @(\"781\u00d78=\", \"821\u00d72=\"),\n    @(\"582\u00d78=\", \"944\u00d75=\"),\n    @(\"193\u00d78=\", \"643\u00d78=\"),\n    @(\"277\u00d75=\", \"815\u00d77=\"),\n    @(\"872\u00d76=\", \"987\u00d79=\"),\n    @(\"419\u00d72=\", \"438\u00d75=\"),\n    @(\"495\u00d76=\", \"262\u00d74=\"),\n    @(\"158\u00d72=\", \"441\u00d74=\"),\n    @(\"152\u00d77=\", \"394\u00d74=\"),\n    @(\"473\u00d79=\", \"997\u00d79=\"),\n    @(\"887\u00d72=\", \"790\u00d76=\"),\n    @(\"400\u00d75=\", \"757\u00d77=\"),\n    @(\"586\u00d73=\", \"681\u00d78=\"),\n    @(\"645\u00d72=\", \"690\u00d73=\"),\n    @(\"584\u00d73=\", \"321\u00d75=\"),\n    @(\"813\u00d79=\", \"822\u00d76=\"),\n    @(\"106\u00d77=\", \"122\u00d72=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
